# Add a new wishlist entry (Astro Bot PS5) as row 38 on the "Valentin" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns are A=Name, B=Image, C=Link, D=Price.
# Set Link/Price first, then Name/Image, so new shared-string entries are
# appended in the same order as the source edit (Link, Price, Name, Image).
$ws.Range("C38").Value = "https://www.wog.ch/index.cfm/details/product/196112%2DAstro%2DBot"
$ws.Range("D38").Value = "64.90 CHF"
$ws.Range("A38").Value = "Astro Bot PS5"
$ws.Range("B38").Value = "https://www.wog.ch/nas/cover_large/p5/p5_astrobot.jpg"

# Reflect the author's next click (cell C39) as the active selection.
$null = $ws.Range("C39").Select()
